# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    populated with the per-fund holding detail for the new quarter (same
#    layout as the other quarterly sheets, e.g. "2021-Q4").
# 2. Insert a new leading row into the "总计" sheet summarising 2022-Q1
#    (date / holding count / holding value), pushing the older rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet, placed immediately before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Match page margins used by every other sheet in the workbook (inches,
# expressed here in points since PageSetup takes points).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Data rows: index column (A) + 4 fund rows
$q1Data = @(
    @("012010", "富国泰享回报6个月持有期混合型证券投资基金A", "9.29", "29.91", "0.67", "0.0622", 10),
    @("001942", "前海开源沪港深汇鑫灵活配置混合A", "0.10", "90.39", "3.16", "0.0032", 10),
    @("001943", "前海开源沪港深汇鑫灵活配置混合C", "0.08", "90.39", "3.16", "0.0025", 10),
    @("012011", "富国泰享回报6个月持有期混合型证券投资基金C", "0.09", "29.91", "0.67", "0.0006", 10)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Range("A$r").Value = ($r - 2)

    # Columns B-G carry their numeric-looking values as text (matches the
    # source data, which keeps fixed decimal formatting like "0.10").
    $q1.Range("B$r`:G$r").NumberFormat = "@"
    $q1.Range("B$r").Value = $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = $row[2]
    $q1.Range("E$r").Value = $row[3]
    $q1.Range("F$r").Value = $row[4]
    $q1.Range("G$r").Value = $row[5]
    $q1.Range("B$r`:G$r").Style = "Normal"

    # Column H (rank) is a genuine number.
    $q1.Range("H$r").Value = $row[6]

    $r = $r + 1
}

# Give the header row the same style as the other quarterly sheets' header.
# (Re-fetch "2021-Q4" fresh by name; sheet handles bound via Worksheets.Item
# track *position*, not identity, and positions shift as sheets are added.)
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Give the index column (A2:A5) the same style as the index column used
# elsewhere in the workbook.
$src2 = $wb.Worksheets.Item("2021-Q4")
$src2.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A6").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.07000000000000001
$totalSheet.Range("B2:D2").Style = "Normal"

# The pre-existing rows' running index (column A) shifts up by one since a
# new row was inserted above them.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Restore the originally-active tab (sheet creation shifts the selection to
# the new sheet as a side effect).
$wb.Worksheets.Item("2020-Q4").Activate()
